$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap mismatched row pairs (F:V) back to correct order ---
$swapPairs = @(
    @(6,7),
    @(15,16),
    @(40,41),
    @(50,51),
    @(53,54),
    @(60,61),
    @(84,85),
    @(89,90),
    @(111,112),
    @(121,122)
)
foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $vals1 = $ws.Range("F$r1" + ":V$r1").Value2
    $vals2 = $ws.Range("F$r2" + ":V$r2").Value2
    $ws.Range("F$r1" + ":V$r1").Value2 = $vals2
    $ws.Range("F$r2" + ":V$r2").Value2 = $vals1
}

# --- Append new match rows 125-138, copying formatting from row 124 ---
$ws.Range("A124:V124").Copy()
$ws.Range("A125:V138").PasteSpecial(-4122)

$row = New-Object 'object[,]' 1,22
$row[0,0] = 124
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45263.47916666666
$row[0,5] = 'Manisa FK'
$row[0,6] = 3
$row[0,7] = 'Altay'
$row[0,8] = 0
$row[0,9] = 1.41
$row[0,10] = '26/11/2023 14:13'
$row[0,11] = 1.32
$row[0,12] = '03/12/2023 11:22'
$row[0,13] = 4.64
$row[0,14] = '26/11/2023 14:13'
$row[0,15] = 5.34
$row[0,16] = '03/12/2023 11:23'
$row[0,17] = 7.28
$row[0,18] = '26/11/2023 14:13'
$row[0,19] = 9.59
$row[0,20] = '03/12/2023 11:23'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-altay/4hVF4ReQ/'
$ws.Range("A125:V125").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 125
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45263.58333333334
$row[0,5] = 'Erzurumspor'
$row[0,6] = 1
$row[0,7] = 'Genclerbirligi'
$row[0,8] = 1
$row[0,9] = 2.83
$row[0,10] = '27/11/2023 18:12'
$row[0,11] = 2.71
$row[0,12] = '03/12/2023 13:59'
$row[0,13] = 3.22
$row[0,14] = '27/11/2023 18:12'
$row[0,15] = 3.19
$row[0,16] = '03/12/2023 13:57'
$row[0,17] = 2.57
$row[0,18] = '27/11/2023 18:12'
$row[0,19] = 2.77
$row[0,20] = '03/12/2023 13:59'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/erzurumspor-fk-genclerbirligi/rws76mRD/'
$ws.Range("A126:V126").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 126
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45263.70833333334
$row[0,5] = 'Goztepe'
$row[0,6] = 3
$row[0,7] = 'Giresunspor'
$row[0,8] = 0
$row[0,9] = 1.34
$row[0,10] = '26/11/2023 17:12'
$row[0,11] = 1.23
$row[0,12] = '03/12/2023 16:55'
$row[0,13] = 4.92
$row[0,14] = '26/11/2023 17:12'
$row[0,15] = 6.07
$row[0,16] = '03/12/2023 16:55'
$row[0,17] = 9.199999999999999
$row[0,18] = '26/11/2023 17:12'
$row[0,19] = 13.22
$row[0,20] = '03/12/2023 16:55'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/goztepe-giresunspor/MTuB57tK/'
$ws.Range("A127:V127").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 127
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45269.47916666666
$row[0,5] = 'Bandirmaspor'
$row[0,6] = 0
$row[0,7] = 'Kocaelispor'
$row[0,8] = 3
$row[0,9] = 2.1
$row[0,10] = '03/12/2023 15:42'
$row[0,11] = 2.12
$row[0,12] = '09/12/2023 11:25'
$row[0,13] = 3.36
$row[0,14] = '03/12/2023 15:42'
$row[0,15] = 3.42
$row[0,16] = '09/12/2023 11:25'
$row[0,17] = 3.34
$row[0,18] = '03/12/2023 15:42'
$row[0,19] = 3.56
$row[0,20] = '09/12/2023 11:25'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-kocaelispor/lbS1eNIC/'
$ws.Range("A128:V128").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 128
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45269.47916666666
$row[0,5] = 'Bodrumspor'
$row[0,6] = 1
$row[0,7] = 'Tuzlaspor'
$row[0,8] = 1
$row[0,9] = 1.36
$row[0,10] = '03/12/2023 15:42'
$row[0,11] = 1.38
$row[0,12] = '09/12/2023 11:26'
$row[0,13] = 4.58
$row[0,14] = '03/12/2023 15:42'
$row[0,15] = 4.73
$row[0,16] = '09/12/2023 11:29'
$row[0,17] = 7.44
$row[0,18] = '03/12/2023 15:42'
$row[0,19] = 8.84
$row[0,20] = '09/12/2023 11:29'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-tuzlaspor/GpUgcql0/'
$ws.Range("A129:V129").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 129
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45269.58333333334
$row[0,5] = 'Umraniyespor'
$row[0,6] = 2
$row[0,7] = 'Keciorengucu'
$row[0,8] = 0
$row[0,9] = 2.08
$row[0,10] = '03/12/2023 15:42'
$row[0,11] = 2.1
$row[0,12] = '09/12/2023 13:57'
$row[0,13] = 3.4
$row[0,14] = '03/12/2023 15:42'
$row[0,15] = 3.51
$row[0,16] = '09/12/2023 13:58'
$row[0,17] = 3.36
$row[0,18] = '03/12/2023 15:42'
$row[0,19] = 3.52
$row[0,20] = '09/12/2023 13:58'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-keciorengucu/b9ZlbPYg/'
$ws.Range("A130:V130").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 130
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45269.70833333334
$row[0,5] = 'Sakaryaspor'
$row[0,6] = 3
$row[0,7] = 'Giresunspor'
$row[0,8] = 1
$row[0,9] = 1.4
$row[0,10] = '03/12/2023 17:13'
$row[0,11] = 1.29
$row[0,12] = '09/12/2023 16:49'
$row[0,13] = 4.45
$row[0,14] = '03/12/2023 17:13'
$row[0,15] = 5.6
$row[0,16] = '09/12/2023 16:56'
$row[0,17] = 6.78
$row[0,18] = '03/12/2023 17:13'
$row[0,19] = 10.16
$row[0,20] = '09/12/2023 16:56'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/sakaryaspor-giresunspor/KIzoa5Jm/'
$ws.Range("A131:V131").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 131
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45270.47916666666
$row[0,5] = 'Corum'
$row[0,6] = 4
$row[0,7] = 'Erzurumspor'
$row[0,8] = 1
$row[0,9] = 1.71
$row[0,10] = '03/12/2023 15:42'
$row[0,11] = 1.91
$row[0,12] = '10/12/2023 08:28'
$row[0,13] = 3.66
$row[0,14] = '03/12/2023 15:42'
$row[0,15] = 3.41
$row[0,16] = '10/12/2023 11:07'
$row[0,17] = 4.54
$row[0,18] = '03/12/2023 15:42'
$row[0,19] = 4.35
$row[0,20] = '10/12/2023 11:07'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/corum-fk-erzurumspor-fk/fyTcd336/'
$ws.Range("A132:V132").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 132
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45270.47916666666
$row[0,5] = 'Manisa FK'
$row[0,6] = 1
$row[0,7] = 'Adanaspor AS'
$row[0,8] = 1
$row[0,9] = 1.63
$row[0,10] = '03/12/2023 15:42'
$row[0,11] = 1.43
$row[0,12] = '10/12/2023 11:19'
$row[0,13] = 3.88
$row[0,14] = '03/12/2023 15:42'
$row[0,15] = 4.63
$row[0,16] = '10/12/2023 11:25'
$row[0,17] = 4.79
$row[0,18] = '03/12/2023 15:42'
$row[0,19] = 7.46
$row[0,20] = '10/12/2023 11:25'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-adanaspor-as/rc9hIpQP/'
$ws.Range("A133:V133").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 133
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45270.58333333334
$row[0,5] = 'Genclerbirligi'
$row[0,6] = 0
$row[0,7] = 'Goztepe'
$row[0,8] = 3
$row[0,9] = 2.66
$row[0,10] = '03/12/2023 18:12'
$row[0,11] = 3.07
$row[0,12] = '10/12/2023 13:58'
$row[0,13] = 3.03
$row[0,14] = '03/12/2023 18:12'
$row[0,15] = 3.21
$row[0,16] = '10/12/2023 13:58'
$row[0,17] = 2.72
$row[0,18] = '03/12/2023 18:12'
$row[0,19] = 2.46
$row[0,20] = '10/12/2023 13:51'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/genclerbirligi-goztepe/Qeys0o4s/'
$ws.Range("A134:V134").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 134
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45270.70833333334
$row[0,5] = 'Eyupspor'
$row[0,6] = 4
$row[0,7] = 'Sanliurfaspor'
$row[0,8] = 0
$row[0,9] = 1.16
$row[0,10] = '03/12/2023 18:12'
$row[0,11] = 1.15
$row[0,12] = '10/12/2023 16:02'
$row[0,13] = 6.66
$row[0,14] = '03/12/2023 18:12'
$row[0,15] = 8.09
$row[0,16] = '10/12/2023 16:33'
$row[0,17] = 12.25
$row[0,18] = '03/12/2023 18:12'
$row[0,19] = 17.71
$row[0,20] = '10/12/2023 16:33'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/eyupspor-sanliurfaspor/21W5fsYI/'
$ws.Range("A135:V135").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 135
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45271.75
$row[0,5] = 'Altay'
$row[0,6] = 0
$row[0,7] = 'Boluspor'
$row[0,8] = 1
$row[0,9] = 3.48
$row[0,10] = '04/12/2023 18:12'
$row[0,11] = 4.99
$row[0,12] = '11/12/2023 17:56'
$row[0,13] = 3.4
$row[0,14] = '04/12/2023 18:12'
$row[0,15] = 3.72
$row[0,16] = '11/12/2023 17:56'
$row[0,17] = 2.04
$row[0,18] = '04/12/2023 18:12'
$row[0,19] = 1.73
$row[0,20] = '11/12/2023 17:56'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/altay-boluspor/46AlJQAJ/'
$ws.Range("A136:V136").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 136
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45279.5
$row[0,5] = 'Tuzlaspor'
$row[0,6] = 1
$row[0,7] = 'Sakaryaspor'
$row[0,8] = 3
$row[0,9] = 3.77
$row[0,10] = '10/12/2023 11:43'
$row[0,11] = 3.45
$row[0,12] = '19/12/2023 11:59'
$row[0,13] = 3.45
$row[0,14] = '10/12/2023 11:43'
$row[0,15] = 3.43
$row[0,16] = '19/12/2023 11:59'
$row[0,17] = 1.92
$row[0,18] = '10/12/2023 11:43'
$row[0,19] = 2.16
$row[0,20] = '19/12/2023 11:59'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-sakaryaspor/MBnu9tmC/'
$ws.Range("A137:V137").Value2 = $row

$row = New-Object 'object[,]' 1,22
$row[0,0] = 137
$row[0,1] = 'turkey'
$row[0,2] = '1-lig'
$row[0,3] = '2023-2024'
$row[0,4] = 45279.625
$row[0,5] = 'Giresunspor'
$row[0,6] = 0
$row[0,7] = 'Genclerbirligi'
$row[0,8] = 1
$row[0,9] = 4.11
$row[0,10] = '10/12/2023 14:12'
$row[0,11] = 5.29
$row[0,12] = '19/12/2023 14:58'
$row[0,13] = 3.48
$row[0,14] = '10/12/2023 14:12'
$row[0,15] = 4.03
$row[0,16] = '19/12/2023 14:58'
$row[0,17] = 1.82
$row[0,18] = '10/12/2023 14:12'
$row[0,19] = 1.63
$row[0,20] = '19/12/2023 14:57'
$row[0,21] = 'https://www.betexplorer.com/football/turkey/1-lig/giresunspor-genclerbirligi/WtUVV3fP/'
$ws.Range("A138:V138").Value2 = $row

